$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "24.673.82"
$ws.Range("E2").Value = "  +1.45%  "

$ws.Range("D3").Value = "1.694.52"
$ws.Range("E3").Value = "  +0.08%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.59%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "316.42"
$ws.Range("E5").Value = "  +1.77%  "

$ws.Range("E6").Value = "  +0.55%  "

$ws.Range("E7").Value = "  +1.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4053"
$ws.Range("E8").Value = "  +0.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.489"
$ws.Range("E9").Value = "  +1.63%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.001"
$ws.Range("E10").Value = "  +0.58%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "52.05"
$ws.Range("E11").Value = "  -4.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08880"
$ws.Range("E12").Value = "  +1.79%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.176"
$ws.Range("E13").Value = "  -0.72%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.52"
$ws.Range("E14").Value = "  +2.23%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.146"
$ws.Range("E15").Value = "  +9.49%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001326"
$ws.Range("E16").Value = "  +0.68%  "

$ws.Range("D17").Value = "1.695.12"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "99.98"
$ws.Range("E18").Value = "  -0.13%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.07006"
$ws.Range("E19").Value = "  +0.07%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.67"
$ws.Range("E20").Value = "  +1.81%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.018"
$ws.Range("E21").Value = "  +5.06%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("E22").Value = "  +0.61%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.37"
$ws.Range("E23").Value = "  +1.74%  "

$ws.Range("D24").Value = "24.674.68"
$ws.Range("E24").Value = "  +1.49%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.207"
$ws.Range("E25").Value = "  +6.38%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.345"
$ws.Range("E26").Value = "  +1.63%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.70"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "163.47"
$ws.Range("E28").Value = "  +2.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "136.26"
$ws.Range("E29").Value = "  +3.32%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.167"
$ws.Range("E30").Value = "  +1.41%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.461"
$ws.Range("E31").Value = "  -3.40%  "

$ws.Range("D32").Value = "1.878.61"
$ws.Range("E32").Value = "  +0.38%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.068"
$ws.Range("E33").Value = "  -1.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08592"
$ws.Range("E34").Value = "  -0.87%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "7.182"
$ws.Range("E35").Value = "  -5.18%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "11.46"
$ws.Range("E36").Value = "  +1.80%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.2738"
$ws.Range("E37").Value = "  +2.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.925"
$ws.Range("E38").Value = "  -1.16%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "14.45"
$ws.Range("E39").Value = "  -1.44%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.09169"
$ws.Range("E40").Value = "  +3.26%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02729"
$ws.Range("E41").Value = "  +0.50%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.473"
$ws.Range("E42").Value = "  +0.83%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7665"
$ws.Range("E43").Value = "  +1.18%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "15.97"
$ws.Range("E44").Value = "  +4.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.600"
$ws.Range("E45").Value = "  +6.25%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.7175"
$ws.Range("E46").Value = "  +0.72%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.225"
$ws.Range("E47").Value = "  +2.21%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.001"
$ws.Range("E48").Value = "  +0.64%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.334"
$ws.Range("E49").Value = "  +6.11%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "140.46"
$ws.Range("E50").Value = "  +0.67%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07975"
$ws.Range("E51").Value = "  +0.90%  "
